$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-03-12 Sunday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-03-13 Monday", 2) | Out-Null
$d.Content.Find.Execute("28+23=", $true, $true, $false, $false, $false, $true, 1, $false, "7+6=", 2) | Out-Null
$d.Content.Find.Execute("21+23=", $true, $true, $false, $false, $false, $true, 1, $false, "36+2=", 2) | Out-Null
$d.Content.Find.Execute("81+12=", $true, $true, $false, $false, $false, $true, 1, $false, "63-53=", 2) | Out-Null
$d.Content.Find.Execute("49+18=", $true, $true, $false, $false, $false, $true, 1, $false, "47+6=", 2) | Out-Null
$d.Content.Find.Execute("33+8=", $true, $true, $false, $false, $false, $true, 1, $false, "78+4=", 2) | Out-Null
$d.Content.Find.Execute("96-11=", $true, $true, $false, $false, $false, $true, 1, $false, "16+23=", 2) | Out-Null
$d.Content.Find.Execute("29-15=", $true, $true, $false, $false, $false, $true, 1, $false, "44+13=", 2) | Out-Null
$d.Content.Find.Execute("18+36=", $true, $true, $false, $false, $false, $true, 1, $false, "49-17=", 2) | Out-Null
$d.Content.Find.Execute("26-8=", $true, $true, $false, $false, $false, $true, 1, $false, "38+40=", 2) | Out-Null
$d.Content.Find.Execute("54+4=", $true, $true, $false, $false, $false, $true, 1, $false, "42-21=", 2) | Out-Null
$d.Content.Find.Execute("67+3=", $true, $true, $false, $false, $false, $true, 1, $false, "7+4=", 2) | Out-Null
$d.Content.Find.Execute("32-4=", $true, $true, $false, $false, $false, $true, 1, $false, "85-23=", 2) | Out-Null
$d.Content.Find.Execute("4+94=", $true, $true, $false, $false, $false, $true, 1, $false, "21+1=", 2) | Out-Null
$d.Content.Find.Execute("50-20=", $true, $true, $false, $false, $false, $true, 1, $false, "83+8=", 2) | Out-Null
$d.Content.Find.Execute("79+4=", $true, $true, $false, $false, $false, $true, 1, $false, "17+81=", 2) | Out-Null
$d.Content.Find.Execute("21+21=", $true, $true, $false, $false, $false, $true, 1, $false, "20+30=", 2) | Out-Null
$d.Content.Find.Execute("15+79=", $true, $true, $false, $false, $false, $true, 1, $false, "16+36=", 2) | Out-Null
$d.Content.Find.Execute("24+11=", $true, $true, $false, $false, $false, $true, 1, $false, "71-58=", 2) | Out-Null
$d.Content.Find.Execute("2-0=", $true, $true, $false, $false, $false, $true, 1, $false, "88-50=", 2) | Out-Null
$d.Content.Find.Execute("59-37=", $true, $true, $false, $false, $false, $true, 1, $false, "43+56=", 2) | Out-Null
$d.Content.Find.Execute("29-11=", $true, $true, $false, $false, $false, $true, 1, $false, "13+69=", 2) | Out-Null
$d.Content.Find.Execute("51-44=", $true, $true, $false, $false, $false, $true, 1, $false, "30-23=", 2) | Out-Null
$d.Content.Find.Execute("39+26=", $true, $true, $false, $false, $false, $true, 1, $false, "3+54=", 2) | Out-Null
$d.Content.Find.Execute("79-39=", $true, $true, $false, $false, $false, $true, 1, $false, "60+7=", 2) | Out-Null
$d.Content.Find.Execute("11+19=", $true, $true, $false, $false, $false, $true, 1, $false, "55+12=", 2) | Out-Null
$d.Content.Find.Execute("34+39=", $true, $true, $false, $false, $false, $true, 1, $false, "63-19=", 2) | Out-Null
$d.Content.Find.Execute("21-6=", $true, $true, $false, $false, $false, $true, 1, $false, "31+40=", 2) | Out-Null
$d.Content.Find.Execute("11-3=", $true, $true, $false, $false, $false, $true, 1, $false, "16+53=", 2) | Out-Null
$d.Content.Find.Execute("65-44=", $true, $true, $false, $false, $false, $true, 1, $false, "91+7=", 2) | Out-Null
$d.Content.Find.Execute("73+0=", $true, $true, $false, $false, $false, $true, 1, $false, "41-11=", 2) | Out-Null
$d.Content.Find.Execute("82-64=", $true, $true, $false, $false, $false, $true, 1, $false, "77-11=", 2) | Out-Null
$d.Content.Find.Execute("29+9=", $true, $true, $false, $false, $false, $true, 1, $false, "74-5=", 2) | Out-Null
$d.Content.Find.Execute("87-55=", $true, $true, $false, $false, $false, $true, 1, $false, "90-15=", 2) | Out-Null
$d.Content.Find.Execute("99-92=", $true, $true, $false, $false, $false, $true, 1, $false, "50+47=", 2) | Out-Null
$d.Content.Find.Execute("23+39=", $true, $true, $false, $false, $false, $true, 1, $false, "46-34=", 2) | Out-Null
$d.Content.Find.Execute("30+18=", $true, $true, $false, $false, $false, $true, 1, $false, "1+10=", 2) | Out-Null
$d.Content.Find.Execute("26+42=", $true, $true, $false, $false, $false, $true, 1, $false, "18+67=", 2) | Out-Null
$d.Content.Find.Execute("94-26=", $true, $true, $false, $false, $false, $true, 1, $false, "19+71=", 2) | Out-Null
$d.Content.Find.Execute("90-85=", $true, $true, $false, $false, $false, $true, 1, $false, "11+17=", 2) | Out-Null
$d.Content.Find.Execute("0+87=", $true, $true, $false, $false, $false, $true, 1, $false, "57-39=", 2) | Out-Null
$d.Content.Find.Execute("20+54=", $true, $true, $false, $false, $false, $true, 1, $false, "29+36=", 2) | Out-Null
$d.Content.Find.Execute("24-7=", $true, $true, $false, $false, $false, $true, 1, $false, "89-24=", 2) | Out-Null
$d.Content.Find.Execute("72-69=", $true, $true, $false, $false, $false, $true, 1, $false, "7-6=", 2) | Out-Null
$d.Content.Find.Execute("77-65=", $true, $true, $false, $false, $false, $true, 1, $false, "69-65=", 2) | Out-Null
$d.Content.Find.Execute("97-61=", $true, $true, $false, $false, $false, $true, 1, $false, "63-35=", 2) | Out-Null
$d.Content.Find.Execute("19+65=", $true, $true, $false, $false, $false, $true, 1, $false, "28+38=", 2) | Out-Null
$d.Content.Find.Execute("54+42=", $true, $true, $false, $false, $false, $true, 1, $false, "31-17=", 2) | Out-Null
$d.Content.Find.Execute("59-25=", $true, $true, $false, $false, $false, $true, 1, $false, "41+41=", 2) | Out-Null
$d.Content.Find.Execute("2+33=", $true, $true, $false, $false, $false, $true, 1, $false, "98-74=", 2) | Out-Null
$d.Content.Find.Execute("80-75=", $true, $true, $false, $false, $false, $true, 1, $false, "84-29=", 2) | Out-Null
$d.Content.Find.Execute("6+80=", $true, $true, $false, $false, $false, $true, 1, $false, "21+38=", 2) | Out-Null
$d.Content.Find.Execute("75-2=", $true, $true, $false, $false, $false, $true, 1, $false, "76+8=", 2) | Out-Null
$d.Content.Find.Execute("74-12=", $true, $true, $false, $false, $false, $true, 1, $false, "21+55=", 2) | Out-Null
$d.Content.Find.Execute("7+77=", $true, $true, $false, $false, $false, $true, 1, $false, "65+34=", 2) | Out-Null
$d.Content.Find.Execute("39-9=", $true, $true, $false, $false, $false, $true, 1, $false, "57-49=", 2) | Out-Null
$d.Content.Find.Execute("33-9=", $true, $true, $false, $false, $false, $true, 1, $false, "72+25=", 2) | Out-Null
$d.Content.Find.Execute("7+18=", $true, $true, $false, $false, $false, $true, 1, $false, "76-2=", 2) | Out-Null
$d.Content.Find.Execute("33-10=", $true, $true, $false, $false, $false, $true, 1, $false, "23-20=", 2) | Out-Null
$d.Content.Find.Execute("54-50=", $true, $true, $false, $false, $false, $true, 1, $false, "61+16=", 2) | Out-Null
$d.Content.Find.Execute("79-22=", $true, $true, $false, $false, $false, $true, 1, $false, "20+79=", 2) | Out-Null
$d.Content.Find.Execute("4+28=", $true, $true, $false, $false, $false, $true, 1, $false, "42+35=", 2) | Out-Null
$d.Content.Find.Execute("87+9=", $true, $true, $false, $false, $false, $true, 1, $false, "13+20=", 2) | Out-Null
$d.Content.Find.Execute("87-0=", $true, $true, $false, $false, $false, $true, 1, $false, "89-30=", 2) | Out-Null
$d.Content.Find.Execute("54+23=", $true, $true, $false, $false, $false, $true, 1, $false, "69-31=", 2) | Out-Null
$d.Content.Find.Execute("58-31=", $true, $true, $false, $false, $false, $true, 1, $false, "45-41=", 2) | Out-Null
$d.Content.Find.Execute("56+6=", $true, $true, $false, $false, $false, $true, 1, $false, "44+42=", 2) | Out-Null
$d.Content.Find.Execute("63-4=", $true, $true, $false, $false, $false, $true, 1, $false, "98-43=", 2) | Out-Null
$d.Content.Find.Execute("65-19=", $true, $true, $false, $false, $false, $true, 1, $false, "6+54=", 2) | Out-Null
$d.Content.Find.Execute("49-46=", $true, $true, $false, $false, $false, $true, 1, $false, "38+22=", 2) | Out-Null
$d.Content.Find.Execute("47+19=", $true, $true, $false, $false, $false, $true, 1, $false, "89-47=", 2) | Out-Null
$d.Content.Find.Execute("35-1=", $true, $true, $false, $false, $false, $true, 1, $false, "29+4=", 2) | Out-Null
$d.Content.Find.Execute("94-30=", $true, $true, $false, $false, $false, $true, 1, $false, "76+1=", 2) | Out-Null
$d.Content.Find.Execute("71-19=", $true, $true, $false, $false, $false, $true, 1, $false, "53-48=", 2) | Out-Null
$d.Content.Find.Execute("91-26=", $true, $true, $false, $false, $false, $true, 1, $false, "15+25=", 2) | Out-Null
$d.Content.Find.Execute("53-21=", $true, $true, $false, $false, $false, $true, 1, $false, "0+37=", 2) | Out-Null
$d.Content.Find.Execute("0+41=", $true, $true, $false, $false, $false, $true, 1, $false, "88-27=", 2) | Out-Null
$d.Content.Find.Execute("62+37=", $true, $true, $false, $false, $false, $true, 1, $false, "39+39=", 2) | Out-Null
$d.Content.Find.Execute("83-53=", $true, $true, $false, $false, $false, $true, 1, $false, "93-93=", 2) | Out-Null
$d.Content.Find.Execute("25+32=", $true, $true, $false, $false, $false, $true, 1, $false, "48-33=", 2) | Out-Null
$d.Content.Find.Execute("51+39=", $true, $true, $false, $false, $false, $true, 1, $false, "70+19=", 2) | Out-Null
$d.Content.Find.Execute("25+64=", $true, $true, $false, $false, $false, $true, 1, $false, "58+18=", 2) | Out-Null
$d.Content.Find.Execute("91-35=", $true, $true, $false, $false, $false, $true, 1, $false, "31-5=", 2) | Out-Null
$d.Content.Find.Execute("95-93=", $true, $true, $false, $false, $false, $true, 1, $false, "33-25=", 2) | Out-Null
$d.Content.Find.Execute("93-62=", $true, $true, $false, $false, $false, $true, 1, $false, "14+82=", 2) | Out-Null
$d.Content.Find.Execute("43+30=", $true, $true, $false, $false, $false, $true, 1, $false, "11+43=", 2) | Out-Null
$d.Content.Find.Execute("89-21=", $true, $true, $false, $false, $false, $true, 1, $false, "90-14=", 2) | Out-Null
$d.Content.Find.Execute("16+39=", $true, $true, $false, $false, $false, $true, 1, $false, "52-43=", 2) | Out-Null
$d.Content.Find.Execute("33+6=", $true, $true, $false, $false, $false, $true, 1, $false, "62-30=", 2) | Out-Null
$d.Content.Find.Execute("88-24=", $true, $true, $false, $false, $false, $true, 1, $false, "82+5=", 2) | Out-Null
$d.Content.Find.Execute("57-9=", $true, $true, $false, $false, $false, $true, 1, $false, "51-14=", 2) | Out-Null
$d.Content.Find.Execute("65-51=", $true, $true, $false, $false, $false, $true, 1, $false, "44+41=", 2) | Out-Null
$d.Content.Find.Execute("12-11=", $true, $true, $false, $false, $false, $true, 1, $false, "16+69=", 2) | Out-Null
$d.Content.Find.Execute("64+4=", $true, $true, $false, $false, $false, $true, 1, $false, "5-2=", 2) | Out-Null
$d.Content.Find.Execute("35+54=", $true, $true, $false, $false, $false, $true, 1, $false, "54-26=", 2) | Out-Null
$d.Content.Find.Execute("4+3=", $true, $true, $false, $false, $false, $true, 1, $false, "8+13=", 2) | Out-Null
$d.Content.Find.Execute("30+39=", $true, $true, $false, $false, $false, $true, 1, $false, "54-24=", 2) | Out-Null
$d.Content.Find.Execute("22+27=", $true, $true, $false, $false, $false, $true, 1, $false, "87-32=", 2) | Out-Null
$d.Content.Find.Execute("68+28=", $true, $true, $false, $false, $false, $true, 1, $false, "33+10=", 2) | Out-Null
$d.Content.Find.Execute("63-21=", $true, $true, $false, $false, $false, $true, 1, $false, "96-63=", 2) | Out-Null
$d.Content.Find.Execute("58+22=", $true, $true, $false, $false, $false, $true, 1, $false, "33+18=", 2) | Out-Null
